$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Foglio1")
$ws1.Activate()
Write-Host ("Type: " + $excel.ActiveWindow.GetType())
$excel.ActiveWindow.ScrollColumn = 4
Write-Host ("ScrollColumn readback=" + $excel.ActiveWindow.ScrollColumn)
$excel.ActiveWindow.ScrollRow = 1
Write-Host ("ScrollRow readback=" + $excel.ActiveWindow.ScrollRow)
$ws1.Range("D1").Select()
